$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bump the "order" column (H) by 1 for the existing top-level menu rows
# that come after the new "Servicios editar" item (rows 10,20,23,28,30,36,41)
$rowsToBump = @(10, 20, 23, 28, 30, 36, 41)
foreach ($r in $rowsToBump) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = $cell.Value2 + 1
}

# Add new row 44: "Servicios editar" submenu under "Servicios" (id 28)
$ws.Cells.Item(44, 1).Value = 44
$ws.Cells.Item(44, 2).Value = 28
$ws.Cells.Item(44, 4).Value = "servicios/editar"
$ws.Cells.Item(44, 9).Value = "Ruta para editar estado de servicio"
$ws.Cells.Item(44, 3).Value = "Servicios editar"
$ws.Cells.Item(44, 5).Value = "minimize"
$ws.Cells.Item(44, 6).Value = "oculto"
$ws.Cells.Item(44, 7).Value = "Digitador"
$ws.Cells.Item(44, 8).Value = 0

# Add new row 45: "Crear respaldo" top-level menu item
$ws.Cells.Item(45, 1).Value = 45
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = "Crear respaldo"
$ws.Cells.Item(45, 4).Value = "crear-respaldo"
$ws.Cells.Item(45, 5).Value = "backup"
$ws.Cells.Item(45, 6).Value = "visible"
$ws.Cells.Item(45, 7).Value = "Administrador"
$ws.Cells.Item(45, 8).Value = 2
$ws.Cells.Item(45, 9).Value = "Menú para crear respaldo"

# Update sheet view to reflect scrolled position / selection seen after the edit
$excel.ActiveWindow.ScrollRow = 37
$ws.Range("C44").Select()
